$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Merge "First Name" (G) + "Last Name" (H) into a single "Full Name" column (G) ---
for ($r = 2; $r -le 7; $r++) {
    $first = $ws.Cells.Item($r, 7).Value()
    $last  = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 7).Value = "$first $last"
}
$ws.Range("G1").Value = "Full Name *"
$ws.Range("A1").Value = "Quantity *"

# --- 2. Snapshot the hyperlink-cell format before mutating (scratch cell well away from
#         any column/row that will shift later on) ---
$ws.Range("F6").Copy()
$ws.Range("A100").PasteSpecial(-4122)  # xlPasteFormats

# --- 3. Delete the now-redundant "Last Name" column (old column H); everything to its
#         right (Pan, Seller Signatory Emails, Bank Account, ...) shifts one column left. ---
$ws.Columns("H").Delete()

# --- 4. The Hyperlinks collection doesn't auto-shift its ref addresses on a column delete,
#         so rebuild it from scratch with the correct (post-shift) addresses. ---
$ws.Hyperlinks.Delete()
$mailAddrs = @(
    "mailto:emp1@investor1.com",
    "mailto:emp1@investor2.com",
    "mailto:emp3@myfirm.com",
    "mailto:emp1@investor1.com",
    "mailto:emp1@investor2.com",
    "mailto:emp1@investor3.com",
    "mailto:emp1@investor4.com",
    "mailto:emp1@investor5.com",
    "mailto:emp1@investor6.com"
)
$cellRefs = @("F6", "F7", "J4", "J6", "J7", "F2", "F3", "F4", "F5")
for ($i = 0; $i -lt $cellRefs.Length; $i++) {
    $ws.Hyperlinks.Add($ws.Range($cellRefs[$i]), $mailAddrs[$i])
}

# --- 5. Hyperlinks.Add() re-stamps a duplicate cell style on every cell it touches;
#         repaint the original hyperlink format (captured in step 2) back over all 9
#         linked cells so they keep looking exactly as they did before. ---
$ws.Range("A100").Copy()
foreach ($ref in $cellRefs) {
    $ws.Range($ref).PasteSpecial(-4122)
}
$ws.Range("A100").Clear()
